$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores its values as plain text in this
# workbook (e.g. "34.130.71", "1.00") because it uses "." as a
# thousands separator and preserves trailing zeros, which are not
# valid/representable as real Excel numbers. Force Text number
# format on each Price cell we touch, right before assigning its
# new value, so the literal text is preserved exactly instead of
# Excel silently converting it to a number.
$ws.Range('D2,D3,D5,D6,D7,D8,D9,D10,D11,D12,D13,D14,D15,D16,D17,D18,D19,D20,D21,D23,D25,D26,D27,D28,D30,D31,D33,D34,D35,D36,D37,D38,D39,D40,D41,D44,D46,D48,D49,D50,D51').NumberFormat = '@'

$ws.Range('D2').Value = '34.062.07'
$ws.Range('E2').Value = '  -0.52%  '
$ws.Range('D3').Value = '1.820.39'
$ws.Range('E3').Value = '  +2.05%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '224.73'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').Value = '0.554'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '31.80'
$ws.Range('E8').Value = '  -0.70%  '
$ws.Range('D9').Value = '0.289'
$ws.Range('E9').Value = '  +3.10%  '
$ws.Range('D10').Value = '0.0720'
$ws.Range('E10').Value = '  +9.37%  '
$ws.Range('D11').Value = '0.0929'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '2.088.21'
$ws.Range('E12').Value = '  +2.36%  '
$ws.Range('D13').Value = '1.839.50'
$ws.Range('E13').Value = '  +3.29%  '
$ws.Range('D14').Value = '10.78'
$ws.Range('E14').Value = '  -3.63%  '
$ws.Range('D15').Value = '0.641'
$ws.Range('E15').Value = '  +2.31%  '
$ws.Range('D16').Value = '34.134.49'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('D17').Value = '4.31'
$ws.Range('E17').Value = '  +2.46%  '
$ws.Range('D18').Value = '69.43'
$ws.Range('E18').Value = '  +0.86%  '
$ws.Range('D19').Value = '249.78'
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('D20').Value = '0.0₃0787'
$ws.Range('E20').Value = '  +6.33%  '
$ws.Range('D21').Value = '11.06'
$ws.Range('E21').Value = '  +6.82%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').Value = '4.25'
$ws.Range('E23').Value = '  +1.30%  '
$ws.Range('E24').Value = '  +0.88%  '
$ws.Range('D25').Value = '160.30'
$ws.Range('E25').Value = '  +2.08%  '
$ws.Range('D26').Value = '16.52'
$ws.Range('E26').Value = '  +0.74%  '
$ws.Range('D27').Value = '7.20'
$ws.Range('E27').Value = '  +2.82%  '
$ws.Range('D28').Value = '0.114'
$ws.Range('E28').Value = '  +0.50%  '
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').Value = '0.0530'
$ws.Range('E30').Value = '  +3.23%  '
$ws.Range('D31').Value = '3.76'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +1.87%  '
$ws.Range('D33').Value = '3.57'
$ws.Range('E33').Value = '  -0.58%  '
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('D35').Value = '1.432.13'
$ws.Range('E35').Value = '  -0.62%  '
$ws.Range('D36').Value = '1.06'
$ws.Range('E36').Value = '  +0.72%  '
$ws.Range('D37').Value = '0.637'
$ws.Range('E37').Value = '  +2.00%  '
$ws.Range('D38').Value = '0.0189'
$ws.Range('E38').Value = '  +0.93%  '
$ws.Range('D39').Value = '0.955'
$ws.Range('E39').Value = '  +7.52%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '2.80'
$ws.Range('E40').Value = '  -1.87%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '81.23'
$ws.Range('E41').Value = '  -2.01%  '
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('E43').Value = '  +4.39%  '
$ws.Range('D44').Value = '6.06'
$ws.Range('E44').Value = '  +3.87%  '
$ws.Range('E45').Value = '  +0.67%  '
$ws.Range('D46').Value = '1.980.25'
$ws.Range('E46').Value = '  +2.00%  '
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('D48').Value = '106.70'
$ws.Range('E48').Value = '  +8.46%  '
$ws.Range('D49').Value = '0.998'
$ws.Range('E49').Value = '  -0.12%  '
$ws.Range('D50').Value = '11.83'
$ws.Range('E50').Value = '  -3.37%  '
$ws.Range('D51').Value = '0.0₆0123'
$ws.Range('E51').Value = '  +4.19%  '
